$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.974.15"
$ws.Range("E2").Value = "  +5.20%  "
$ws.Range("D3").Value = "2.349.83"
$ws.Range("E3").Value = "  +4.56%  "
$ws.Range("E4").Value = "  -0.74%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.79"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.10"
$ws.Range("E6").Value = "  +3.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.577"
$ws.Range("E7").Value = "  +1.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.539"
$ws.Range("E9").Value = "  +4.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.80"
$ws.Range("E10").Value = "  +2.39%  "
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.44"
$ws.Range("E12").Value = "  +3.11%  "
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").Value = "2.706.62"
$ws.Range("D15").Value = "2.344.68"
$ws.Range("E15").Value = "  +4.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.22"
$ws.Range("E16").Value = "  +4.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.832"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").Value = "46.809.52"
$ws.Range("E18").Value = "  +5.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.76"
$ws.Range("E19").Value = "  +17.28%  "
$ws.Range("E20").Value = "  +1.39%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.84"
$ws.Range("E22").Value = "  +2.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "245.07"
$ws.Range("E23").Value = "  +3.11%  "
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "41.77"
$ws.Range("E27").Value = "  +12.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.30"
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.89"
$ws.Range("E29").Value = "  +1.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.17"
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.77"
$ws.Range("E31").Value = "  -2.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.45"
$ws.Range("E32").Value = "  +3.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0818"
$ws.Range("E33").Value = "  +4.27%  "
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("E35").Value = "  -1.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.109"
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("E38").Value = "  -1.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.04"
$ws.Range("E39").Value = "  +7.10%  "
$ws.Range("E40").Value = "  +5.23%  "
$ws.Range("E41").Value = "  +2.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.81"
$ws.Range("E42").Value = "  -9.32%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.99"
$ws.Range("E44").Value = "  +12.90%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.863.68"
$ws.Range("E45").Value = "  +2.87%  "
$ws.Range("E46").Value = "  +5.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "74.36"
$ws.Range("E47").Value = "  +7.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "81.12"
$ws.Range("E48").Value = "  -0.84%  "
$ws.Range("E49").Value = "  +2.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "98.69"
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.42"
$ws.Range("E51").Value = "  +2.44%  "
